$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Entregas")
$ws2.Columns.Item(4).ColumnWidth = 19.72
Write-Host "OK"
